# TC02_Trials_Filter_Ethnicity-NotHispLatino.xlsx
# "all single filter scripts in CTDC" - update CasesTab queries to the new
# file-centric Cypher pattern, and add a new FilesTab row with file list /
# stat queries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (CasesTab): refresh the query (B2) and stat query (C2) text ---
$ws.Range("B2").Value = "MATCH (c:case)`n MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)`n MATCH (f:file)-[*]->(c)`n  WHERE c.ethnicity =`"NOT_HISPANIC_OR_LATINO`" `nRETURN DISTINCT`n    c.case_id AS ``Case ID``,`n     ct.clinical_trial_designation AS ``Trial Code``,`n     a.arm_id AS Arm,`n      a.arm_drug AS ``Arm Treatment``,`nc.disease AS Diagnosis,`n  c.gender AS Gender,`n    c.race AS Race,`n    c.ethnicity AS Ethnicity"
$ws.Range("C2").Value = "MATCH (f:file)`nOPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)`nOPTIONAL MATCH (f)-[*]->(c:case)`nWITH f,a,ct,c`n      WHERE c.ethnicity =`"NOT_HISPANIC_OR_LATINO`" `nRETURN`n    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,`n    COUNT(DISTINCT c.case_id) AS Cases,`n    COUNT(DISTINCT f) AS Files"

# --- Row 3 (new FilesTab row) ---
$ws.Range("A3").Value = "FilesTab"
$ws.Range("B3").Value = "MATCH (f:file)`nOPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)`nOPTIONAL MATCH (f)-[*]->(c:case)`nOPTIONAL MATCH (f)-->(parent)`nWITH f,a,ct,c,parent`nWHERE c.ethnicity =`"NOT_HISPANIC_OR_LATINO`" `nWITH`n    f, parent, c, a, ct,`n    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n    toInteger(floor(log(f.file_size)/log(1024))) as i,`n    2 as precision`nWITH`n    f, parent, c, a, ct,`n    f.file_size /(1024^i) AS value,`n    10^precision AS factor,`n    units[i] as unit`nWITH`n    f, parent, c, a, ct, unit,`n    round(factor * value)/factor AS size`nRETURN DISTINCT`n    f.file_name AS ``File Name``,`n    head(labels(parent)) as Association,`n    f.file_description AS Description,`n    f.file_format AS ``File Format``,`n    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n    ct.clinical_trial_designation AS ``Trial Code``,`n    a.arm_id AS Arm,`n    c.case_id AS ``Case ID``"
$ws.Range("C3").Value = "MATCH (f:file)`nOPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)`nOPTIONAL MATCH (f)-[*]->(c:case)`nWITH f,a,ct,c`n        WHERE c.ethnicity =`"NOT_HISPANIC_OR_LATINO`" `nRETURN`n    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,`n    COUNT(DISTINCT c.case_id) AS Cases,`n    COUNT(DISTINCT f) AS Files"
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2

# B3/C3 use the same wrapped-text style as B2/C2
$ws.Range("B3:C3").WrapText = $true

# --- Row heights to fit the long wrapped query text ---
$ws.Rows.Item(2).RowHeight = 188.5
$ws.Rows.Item(3).RowHeight = 409.5

# --- View state: zoom out and select/scroll to the new row ---
$excel.ActiveWindow.Zoom = 55
$ws.Range("C3").Select()
